# Update header labels on the existing sheets.
$wb = $excel.ActiveWorkbook

$weekly = $wb.Worksheets.Item("Weekly Quantity")
$weekly.Range("B1").Value = "Weekly_PO_Qty"

$monthly = $wb.Worksheets.Item("Monthly Trend")
$monthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$forecast = $wb.Worksheets.Add($null, $lastSheet)
$forecast.Name = "PO Forecast"

# Header row.
$forecast.Range("A1").Value = "ds"
$forecast.Range("B1").Value = "PO_Forecast"
$forecast.Range("C1").Value = "yhat_lower"
$forecast.Range("D1").Value = "yhat_upper"

$headerRange = $forecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows.
$data = @(
    @(44983.99999999999, 11, 9.700186189984345, 11.85248807989762),
    @(44997.99999999999, 10, 9.261049761434197, 11.46650510803554),
    @(45004.99999999999, 10, 9.15742709001684, 11.22157592547203),
    @(45011.99999999999, 10, 8.873899761378324, 11.09540140741677),
    @(45032.99999999999, 9, 8.314981867805932, 10.50396730069762),
    @(45074.99999999999, 8, 7.255487597370445, 9.436483759882615),
    @(45158.99999999999, 6, 4.863400026039778, 7.154056947396551),
    @(45165.99999999999, 6, 4.730660658698565, 6.91990369134753),
    @(45172.99999999999, 6, 4.572220948545121, 6.732225666535581),
    @(45179.99999999999, 5, 4.473103626922346, 6.608801449161947),
    @(45186.99999999999, 5, 4.172300209848599, 6.362096153519151),
    @(45193.99999999999, 5, 4.065123971191731, 6.172411784771924),
    @(45200.99999999999, 5, 3.818428606665114, 6.050665862546493),
    @(45207.99999999999, 5, 3.566237166581967, 5.819159347669398),
    @(45214.99999999999, 5, 3.439538926731941, 5.578072501302114)
)

$row = 2
foreach ($r in $data) {
    $forecast.Cells.Item($row, 1).Value = $r[0]
    $forecast.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $forecast.Cells.Item($row, 2).Value = $r[1]
    $forecast.Cells.Item($row, 3).Value = $r[2]
    $forecast.Cells.Item($row, 4).Value = $r[3]
    $row = $row + 1
}

# Restore the originally-active sheet/selection.
$weekly.Activate()
[void]$weekly.Range("A1").Select()

Write-Output "done"
